$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Bottom three rows (44-46): collapse the multi-run/tab cell text into a
# single value each ---
$t.Rows.Item(44).Cells.Item(1).Range.Text = "99.95"
$t.Rows.Item(45).Cells.Item(1).Range.Text = "0.37"
$t.Rows.Item(46).Cells.Item(1).Range.Text = "775"

# --- Remove the row that used to hold "0.04366" (old row 12) ---
$t.Rows.Item(12).Delete()

# --- Simple text edits on rows 1-11 (row 12 already removed, rows below are
# unaffected since they were already re-written above) ---
$t.Rows.Item(1).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(2).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(3).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(4).Cells.Item(1).Range.Text = "2012"
$t.Rows.Item(6).Cells.Item(1).Range.Text = "0.00194"
$t.Rows.Item(8).Cells.Item(1).Range.Text = "0.00026"
$t.Rows.Item(9).Cells.Item(1).Range.Text = "0.00031"
$t.Rows.Item(10).Cells.Item(1).Range.Text = "0.00038"
$t.Rows.Item(11).Cells.Item(1).Range.Text = "0.36840"

# --- Insert a new row right after row 6 (the "0.00194" row) carrying
# "0.00018", matching formatting of the surrounding rows ---
$newRow = $t.Rows.Add($t.Rows.Item(7))
$newRow.Cells.Item(1).Range.Text = "0.00018"
